$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3999.25
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 3999.25
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H87").Value = 118671.71
$ws.Range("J87").Value = 118671.71
$ws.Range("L87").Value = 118671.71
$ws.Range("N87").Value = -121167.71
$ws.Range("H90").Value = 118671.71
$ws.Range("J90").Value = 118671.71
$ws.Range("L90").Value = 356015.13
$ws.Range("N90").Value = -368495.13
$ws.Range("H106").Value = 5694.4287
$ws.Range("I106").Value = 6227
$ws.Range("J106").Value = 2499
$ws.Range("K106").Value = 6227
$ws.Range("L106").Value = 2499
$ws.Range("M106").Value = -5596
$ws.Range("N106").Value = -3761
$ws.Range("H137").Value = 5257.302
$ws.Range("I137").Value = 2023.1
$ws.Range("J137").Value = 17696.54
$ws.Range("K137").Value = 6069.299999999999
$ws.Range("L137").Value = 53089.62
$ws.Range("M137").Value = -3519.299999999999
$ws.Range("N137").Value = -58189.62

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5211616.5
$ws.Range("I32").Value = 7577866
$ws.Range("K32").Value = 7577866
$ws.Range("M32").Value = -7577579
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H61").Value = 963823.6
$ws.Range("I61").Value = 3040
$ws.Range("J61").Value = 2885390.8
$ws.Range("K61").Value = 3040
$ws.Range("L61").Value = 2885390.8
$ws.Range("M61").Value = -2828
$ws.Range("N61").Value = -2885814.8
$ws.Range("H122").Value = 3772.875
$ws.Range("I122").Value = 1480.8572
$ws.Range("K122").Value = 4442.571599999999
$ws.Range("M122").Value = -1992.571599999999
$ws.Range("H136").Value = 963823.6
$ws.Range("I136").Value = 3040
$ws.Range("J136").Value = 2885390.8
$ws.Range("K136").Value = 9120
$ws.Range("L136").Value = 8656172.399999999
$ws.Range("M136").Value = -6570
$ws.Range("N136").Value = -8661272.399999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 13944.037
$ws.Range("I99").Value = 14978.75
$ws.Range("K99").Value = 14978.75
$ws.Range("M99").Value = -13480.75
$ws.Range("H107").Value = 1340.6471
$ws.Range("I107").Value = 1486.3334
$ws.Range("J107").Value = 1261.1818
$ws.Range("K107").Value = 1486.3334
$ws.Range("L107").Value = 1261.1818
$ws.Range("M107").Value = 433.6666
$ws.Range("N107").Value = -5101.1818

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13476.23
$ws.Range("I31").Value = 888.3333
$ws.Range("K31").Value = 888.3333
$ws.Range("M31").Value = -593.3333
$ws.Range("H34").Value = 13476.23
$ws.Range("I34").Value = 888.3333
$ws.Range("K34").Value = 888.3333
$ws.Range("M34").Value = -686.3333
$ws.Range("H99").Value = 3815.2222
$ws.Range("I99").Value = 1967.6
$ws.Range("J99").Value = 6124.75
$ws.Range("K99").Value = 1967.6
$ws.Range("L99").Value = 6124.75
$ws.Range("M99").Value = -469.5999999999999
$ws.Range("N99").Value = -9120.75
$ws.Range("H107").Value = 5995
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = -8080
$ws.Range("N107").Value = -5830
$ws.Range("H126").Value = 3815.2222
$ws.Range("I126").Value = 1967.6
$ws.Range("J126").Value = 6124.75
$ws.Range("K126").Value = 5902.799999999999
$ws.Range("L126").Value = 18374.25
$ws.Range("M126").Value = -3432.799999999999
$ws.Range("N126").Value = -23314.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 2899.4
$ws.Range("I26").Value = 2999.5
$ws.Range("J26").Value = 2499
$ws.Range("K26").Value = 8998.5
$ws.Range("L26").Value = 7497
$ws.Range("M26").Value = -8710.5
$ws.Range("N26").Value = -8073
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("H68").Value = 692.3333
$ws.Range("J68").Value = 1097
$ws.Range("L68").Value = 3291
$ws.Range("N68").Value = -4913
$ws.Range("H71").Value = 692.3333
$ws.Range("J71").Value = 1097
$ws.Range("L71").Value = 9873
$ws.Range("N71").Value = -17985
$ws.Range("H81").Value = 10202603
$ws.Range("J81").Value = 17001666
$ws.Range("L81").Value = 51004998
$ws.Range("N81").Value = -51007244
$ws.Range("H84").Value = 10202603
$ws.Range("J84").Value = 17001666
$ws.Range("L84").Value = 153014994
$ws.Range("N84").Value = -153026226
$ws.Range("H92").Value = 279.33334
$ws.Range("J92").Value = 325
$ws.Range("L92").Value = 975
$ws.Range("N92").Value = -3471
$ws.Range("H109").Value = 4166934.5
$ws.Range("I109").Value = 306.2857
$ws.Range("J109").Value = 33333332
$ws.Range("K109").Value = 918.8571000000001
$ws.Range("L109").Value = 99999996
$ws.Range("M109").Value = 121.1428999999999
$ws.Range("N109").Value = -100002076
$ws.Range("H110").Value = 9458.833000000001
$ws.Range("I110").Value = 438.25
$ws.Range("J110").Value = 27500
$ws.Range("K110").Value = 1314.75
$ws.Range("L110").Value = 82500
$ws.Range("M110").Value = 2775.25
$ws.Range("N110").Value = -90680
$ws.Range("H111").Value = 943.2857
$ws.Range("I111").Value = 943.2857
$ws.Range("K111").Value = 2829.8571
$ws.Range("M111").Value = 237.1428999999998
$ws.Range("H112").Value = 3637.75
$ws.Range("I112").Value = 2350.3333
$ws.Range("K112").Value = 7050.999899999999
$ws.Range("M112").Value = -5942.999899999999
$ws.Range("H122").Value = 13453800
$ws.Range("J122").Value = 2581006
$ws.Range("L122").Value = 23229054
$ws.Range("N122").Value = -23233954

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2670.875
$ws.Range("I122").Value = 2816.111
$ws.Range("K122").Value = 8448.332999999999
$ws.Range("M122").Value = -5998.332999999999
$ws.Range("H123").Value = 55236.625
$ws.Range("J123").Value = 55236.625
$ws.Range("L123").Value = 55236.625
$ws.Range("N123").Value = -60136.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2230.6191
$ws.Range("J82").Value = 1592.5
$ws.Range("L82").Value = 1592.5
$ws.Range("N82").Value = -2314.5
$ws.Range("H85").Value = 2230.6191
$ws.Range("J85").Value = 1592.5
$ws.Range("L85").Value = 1592.5
$ws.Range("N85").Value = -4088.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3379.6155
$ws.Range("I122").Value = 1978.3334
$ws.Range("K122").Value = 5935.0002
$ws.Range("M122").Value = -3485.0002
$ws.Range("H132").Value = 496309.4
$ws.Range("I132").Value = 2653.5334
$ws.Range("K132").Value = 7960.600199999999
$ws.Range("M132").Value = -5430.600199999999
